$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure the price/volume columns retain their original text formatting
# (many values look numeric, e.g. "306.27", and would otherwise be
# auto-converted to actual numbers by Excel).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '46.549.86'
$ws.Range("E2").Value = '  +0.84%  '
$ws.Range("D3").Value = '2.587.65'
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '306.27'
$ws.Range("E5").Value = '  +1.43%  '
$ws.Range("D6").Value = '100.23'
$ws.Range("E6").Value = '  +0.53%  '
$ws.Range("D7").Value = '0.602'
$ws.Range("E7").Value = '  +5.46%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = '0.578'
$ws.Range("E9").Value = '  +12.78%  '
$ws.Range("D10").Value = '38.67'
$ws.Range("E10").Value = '  +12.47%  '
$ws.Range("E11").Value = '  +5.14%  '
$ws.Range("D12").Value = '8.16'
$ws.Range("E12").Value = '  +14.60%  '
$ws.Range("D13").Value = '2.984.20'
$ws.Range("E13").Value = '  +9.96%  '
$ws.Range("E14").Value = '  +1.87%  '
$ws.Range("D15").Value = '2.605.93'
$ws.Range("E15").Value = '  +10.79%  '
$ws.Range("D16").Value = '0.907'
$ws.Range("E16").Value = '  +12.15%  '
$ws.Range("D17").Value = '14.90'
$ws.Range("E17").Value = '  +9.08%  '
$ws.Range("D18").Value = '46.733.90'
$ws.Range("E18").Value = '  +1.48%  '
$ws.Range("D19").Value = '13.39'
$ws.Range("E19").Value = '  +4.78%  '
$ws.Range("E20").Value = '  +4.87%  '
$ws.Range("D21").Value = '6.66'
$ws.Range("E21").Value = '  +10.22%  '
$ws.Range("D22").Value = '71.24'
$ws.Range("E22").Value = '  +5.76%  '
$ws.Range("D23").Value = '255.44'
$ws.Range("E23").Value = '  +3.70%  '
$ws.Range("E24").Value = '  +4.74%  '
$ws.Range("E25").Value = '  +13.95%  '
$ws.Range("D26").Value = '28.26'
$ws.Range("E26").Value = '  +34.92%  '
$ws.Range("E27").Value = '  +0.15%  '
$ws.Range("D28").Value = '10.51'
$ws.Range("E28").Value = '  +7.26%  '
$ws.Range("D29").Value = '2.30'
$ws.Range("E29").Value = '  +4.97%  '
$ws.Range("D30").Value = '39.69'
$ws.Range("E30").Value = '  -0.62%  '
$ws.Range("D31").Value = '3.78'
$ws.Range("E31").Value = '  +2.17%  '
$ws.Range("D32").Value = '6.20'
$ws.Range("E32").Value = '  +11.98%  '
$ws.Range("E33").Value = '  +24.57%  '
$ws.Range("D34").Value = '2.93'
$ws.Range("E34").Value = '  +5.61%  '
$ws.Range("D35").Value = '0.0835'
$ws.Range("E35").Value = '  +7.83%  '
$ws.Range("D36").Value = '150.19'
$ws.Range("E36").Value = '  +2.77%  '
$ws.Range("D37").Value = '0.116'
$ws.Range("E37").Value = '  +3.13%  '
$ws.Range("E38").Value = '  +4.81%  '
$ws.Range("D39").Value = '4.20'
$ws.Range("E39").Value = '  +6.01%  '
$ws.Range("D40").Value = '15.80'
$ws.Range("E40").Value = '  +5.50%  '
$ws.Range("D41").Value = '3.66'
$ws.Range("E41").Value = '  +13.21%  '
$ws.Range("D42").Value = '0.0323'
$ws.Range("E42").Value = '  +7.45%  '
$ws.Range("D43").Value = '2.029.15'
$ws.Range("E43").Value = '  +8.04%  '
$ws.Range("D44").Value = '18.83'
$ws.Range("E44").Value = '  +31.85%  '
$ws.Range("D45").Value = '0.998'
$ws.Range("E45").Value = '  -0.02%  '
$ws.Range("D46").Value = '92.00'
$ws.Range("E46").Value = '  -1.56%  '
$ws.Range("E47").Value = '  +0.38%  '
$ws.Range("D48").Value = '9.12'
$ws.Range("E48").Value = '  +10.49%  '
$ws.Range("D49").Value = '109.15'
$ws.Range("E49").Value = '  +11.79%  '
$ws.Range("E50").Value = '  +8.17%  '
$ws.Range("D51").Value = '2.842.75'
$ws.Range("E51").Value = '  +9.97%  '
